$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.72626133333333
$ws.Range("H2").Value = 38.178784
$ws.Range("I2").Value = 0.08664518826379554
$ws.Range("J2").Value = 0.08664518826379553
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7479136666666667
$ws.Range("N2").Value = 2.243741
$ws.Range("Q2").Value = 9.518144776771557
$ws.Range("R2").Value = 85.663302990944
$ws.Range("S2").Value = 0.08664518826379554
$ws.Range("T2").Value = 0.08664518826379553

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 132.5150096666667
$ws.Range("H3").Value = 397.545029
$ws.Range("I3").Value = 0.9022121783931373
$ws.Range("J3").Value = 0.9022121783931372
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.7479136666666667
$ws.Range("N3").Value = 2.243741
$ws.Range("Q3").Value = 99.10978676816543
$ws.Range("R3").Value = 891.9880809134889
$ws.Range("S3").Value = 0.9022121783931373
$ws.Range("T3").Value = 0.9022121783931372

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.636606333333333
$ws.Range("H4").Value = 4.909819
$ws.Range("I4").Value = 0.0111426333430672
$ws.Range("J4").Value = 0.01114263334306719
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.7479136666666667
$ws.Range("N4").Value = 2.243741
$ws.Range("Q4").Value = 1.224040243653222
$ws.Range("R4").Value = 11.016362192879
$ws.Range("S4").Value = 0.0111426333430672
$ws.Range("T4").Value = 0.01114263334306719
